# Auto-generated edit script: updates market-price-derived profit
# columns (H:N) on several rows across all 8 job sheets to match
# the latest scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 160.3
$ws.Range("I9").Value = 111.55556
$ws.Range("K9").Value = 111.55556
$ws.Range("M9").Value = 57.44444

# Row 40
$ws.Range("H40").Value = 83336664
$ws.Range("I40").Value = 4000.75
$ws.Range("J40").Value = 250001980
$ws.Range("K40").Value = 4000.75
$ws.Range("L40").Value = 250001980
$ws.Range("M40").Value = -3825.75
$ws.Range("N40").Value = -250002330

# Row 58
$ws.Range("H58").Value = 5242.7144
$ws.Range("I58").Value = 841
$ws.Range("J58").Value = 7003.4
$ws.Range("K58").Value = 2523
$ws.Range("L58").Value = 21010.2
$ws.Range("M58").Value = -2373
$ws.Range("N58").Value = -21310.2

# Row 96
$ws.Range("H96").Value = 2904482
$ws.Range("I96").Value = 5385.6665
$ws.Range("J96").Value = 7253126.5
$ws.Range("K96").Value = 16156.9995
$ws.Range("L96").Value = 21759379.5
$ws.Range("M96").Value = -14783.9995
$ws.Range("N96").Value = -21762125.5

# Row 100
$ws.Range("H100").Value = 4454.8335
$ws.Range("I100").Value = 2544.6365
$ws.Range("K100").Value = 2544.6365
$ws.Range("M100").Value = -2003.6365

# Row 124
$ws.Range("H124").Value = 99998.5
$ws.Range("J124").Value = 99998.5
$ws.Range("L124").Value = 99998.5
$ws.Range("N124").Value = -109818.5

# Row 132
$ws.Range("H132").Value = 2554.2666
$ws.Range("I132").Value = 2325.8215
$ws.Range("K132").Value = 6977.4645
$ws.Range("M132").Value = -4447.4645

# Row 137
$ws.Range("H137").Value = 3162.6924
$ws.Range("I137").Value = 3458.1875
$ws.Range("K137").Value = 10374.5625
$ws.Range("M137").Value = -7824.5625


$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5547.3555
$ws.Range("I32").Value = 5650.386
$ws.Range("K32").Value = 5650.386
$ws.Range("M32").Value = -5363.386

# Row 45
$ws.Range("H45").Value = 2781.6191
$ws.Range("I45").Value = 1806.2354
$ws.Range("J45").Value = 6927
$ws.Range("K45").Value = 1806.2354
$ws.Range("L45").Value = 6927
$ws.Range("M45").Value = -1429.2354
$ws.Range("N45").Value = -7681

# Row 110
$ws.Range("H110").Value = 5937.893
$ws.Range("I110").Value = 6308.7617
$ws.Range("K110").Value = 6308.7617
$ws.Range("M110").Value = -4263.7617

# Row 120
$ws.Range("H120").Value = 128000
$ws.Range("J120").Value = 128000
$ws.Range("L120").Value = 128000
$ws.Range("N120").Value = -137676

# Row 122
$ws.Range("H122").Value = 2836.4736
$ws.Range("J122").Value = 2795.1428
$ws.Range("L122").Value = 8385.428400000001
$ws.Range("N122").Value = -13285.4284


$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 22197.8
$ws.Range("I82").Value = 5247.25
$ws.Range("J82").Value = 90000
$ws.Range("K82").Value = 5247.25
$ws.Range("L82").Value = 90000
$ws.Range("M82").Value = -4864.25
$ws.Range("N82").Value = -90766

# Row 85
$ws.Range("H85").Value = 22197.8
$ws.Range("I85").Value = 5247.25
$ws.Range("J85").Value = 90000
$ws.Range("K85").Value = 5247.25
$ws.Range("L85").Value = 90000
$ws.Range("M85").Value = -3921.25
$ws.Range("N85").Value = -92652

# Row 99
$ws.Range("H99").Value = 2500
$ws.Range("I99").Value = 2500
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2500
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1002
$ws.Range("N99").ClearContents()

# Row 103
$ws.Range("H103").Value = 14878
$ws.Range("J103").Value = 14878
$ws.Range("L103").Value = 14878
$ws.Range("N103").Value = -17222

# Row 107
$ws.Range("H107").Value = 8784
$ws.Range("I107").Value = 10180
$ws.Range("K107").Value = 10180
$ws.Range("M107").Value = -8260


$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1729.1
$ws.Range("I16").Value = 983
$ws.Range("K16").Value = 983
$ws.Range("M16").Value = -696

# Row 31
$ws.Range("H31").Value = 20411216
$ws.Range("I31").Value = 30305564
$ws.Range("K31").Value = 30305564
$ws.Range("M31").Value = -30305269

# Row 34
$ws.Range("H34").Value = 20411216
$ws.Range("I34").Value = 30305564
$ws.Range("K34").Value = 30305564
$ws.Range("M34").Value = -30305362

# Row 113
$ws.Range("H113").Value = 1729.1
$ws.Range("I113").Value = 983
$ws.Range("K113").Value = 983
$ws.Range("M113").Value = 1187

# Row 134
$ws.Range("H134").Value = 1970.04
$ws.Range("I134").Value = 1993.25
$ws.Range("J134").Value = 1877.2
$ws.Range("K134").Value = 5979.75
$ws.Range("L134").Value = 5631.6
$ws.Range("M134").Value = -3444.75
$ws.Range("N134").Value = -10701.6


$ws = $wb.Worksheets.Item("CUL")
# Row 18
$ws.Range("H18").Value = 403.9
$ws.Range("I18").Value = 337.66666
$ws.Range("K18").Value = 1012.99998
$ws.Range("M18").Value = -843.9999799999999

# Row 39
$ws.Range("H39").Value = 11337
$ws.Range("J39").Value = 12666.375
$ws.Range("L39").Value = 37999.125
$ws.Range("N39").Value = -38587.125

# Row 60
$ws.Range("H60").Value = 3884.7778
$ws.Range("I60").Value = 203.75
$ws.Range("J60").Value = 33333
$ws.Range("K60").Value = 611.25
$ws.Range("L60").Value = 99999
$ws.Range("M60").Value = -360.25
$ws.Range("N60").Value = -100501

# Row 120
$ws.Range("H120").Value = 30817
$ws.Range("I120").Value = 29968
$ws.Range("K120").Value = 89904
$ws.Range("M120").Value = -85066

# Row 138
$ws.Range("H138").Value = 13126.315
$ws.Range("I138").Value = 13641.111
$ws.Range("K138").Value = 40923.333
$ws.Range("M138").Value = -35783.333


$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3536.125
$ws.Range("I80").Value = 3473.75
$ws.Range("J80").Value = 3598.5
$ws.Range("K80").Value = 3473.75
$ws.Range("L80").Value = 3598.5
$ws.Range("M80").Value = -2475.75
$ws.Range("N80").Value = -5594.5

# Row 83
$ws.Range("H83").Value = 3536.125
$ws.Range("I83").Value = 3473.75
$ws.Range("J83").Value = 3598.5
$ws.Range("K83").Value = 17368.75
$ws.Range("L83").Value = 17992.5
$ws.Range("M83").Value = -12376.75
$ws.Range("N83").Value = -27976.5

# Row 97
$ws.Range("H97").Value = 1288.7142
$ws.Range("I97").Value = 1871.8334
$ws.Range("K97").Value = 1871.8334
$ws.Range("M97").Value = -1375.8334

# Row 102
$ws.Range("H102").Value = 2460.3704
$ws.Range("I102").Value = 2260.9
$ws.Range("K102").Value = 2260.9
$ws.Range("M102").Value = -638.9000000000001

# Row 126
$ws.Range("H126").Value = 3362.3333
$ws.Range("I126").Value = 3058
$ws.Range("K126").Value = 9174
$ws.Range("M126").Value = -6704


$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 9603.723
$ws.Range("I7").Value = 8239
$ws.Range("K7").Value = 8239
$ws.Range("M7").Value = -8127

# Row 46
$ws.Range("H46").Value = 2488.5
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

# Row 63
$ws.Range("H63").Value = 99599.8
$ws.Range("J63").Value = 99599.8
$ws.Range("L63").Value = 99599.8
$ws.Range("N63").Value = -101097.8

# Row 66
$ws.Range("H66").Value = 99599.8
$ws.Range("J66").Value = 99599.8
$ws.Range("L66").Value = 298799.4
$ws.Range("N66").Value = -306287.4

# Row 82
$ws.Range("H82").Value = 5720.6
$ws.Range("I82").Value = 6200.6665
$ws.Range("K82").Value = 6200.6665
$ws.Range("M82").Value = -5839.6665

# Row 85
$ws.Range("H85").Value = 5720.6
$ws.Range("I85").Value = 6200.6665
$ws.Range("K85").Value = 6200.6665
$ws.Range("M85").Value = -4952.6665

# Row 100
$ws.Range("H100").Value = 11919849
$ws.Range("I100").Value = 2777.889
$ws.Range("K100").Value = 2777.889
$ws.Range("M100").Value = -2236.889

# Row 126
$ws.Range("H126").Value = 9603.723
$ws.Range("I126").Value = 8239
$ws.Range("K126").Value = 24717
$ws.Range("M126").Value = -22247

# Row 136
$ws.Range("H136").Value = 4223.9165
$ws.Range("I136").Value = 1684
$ws.Range("J136").Value = 7779.8
$ws.Range("K136").Value = 5052
$ws.Range("L136").Value = 23339.4
$ws.Range("M136").Value = -2502
$ws.Range("N136").Value = -28439.4


$ws = $wb.Worksheets.Item("WVR")
# Row 101
$ws.Range("H101").Value = 35903
$ws.Range("J101").Value = 35903
$ws.Range("L101").Value = 35903
$ws.Range("N101").Value = -42393

# Row 140
$ws.Range("H140").Value = 51426.715
$ws.Range("J140").Value = 54997.832
$ws.Range("L140").Value = 54997.832
$ws.Range("N140").Value = -65357.832

